# Insert a new "Foto" column into the Contenedores inventory sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Contenedores")

# Insert a new column at G, shifting Costo/Stock/Bodega/Ultima Actualizacion right.
$ws.Columns.Item(7).Insert()

# Header for the new column - copy the header formatting from the
# neighboring (already-styled) header cell, then set the text.
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)
$ws.Cells.Item(1, 7).Value = "Foto"

# Photo filenames per row (rows 2-16).
$fotos = @(
    "florero-vidrio-grande.jpg",
    "florero-vidrio-mediano.jpg",
    "florero-vidrio-pequeno.jpg",
    "florero-ceramica-blanco.jpg",
    "florero-ceramica-negro.jpg",
    "macetero-terracota-grande.jpg",
    "macetero-terracota-mediano.jpg",
    "macetero-terracota-pequeno.jpg",
    "macetero-plastico-blanco.jpg",
    "canasto-rectangular-grande.jpg",
    "canasto-rectangular-mediano.jpg",
    "canasto-redondo-pequeno.jpg",
    "florero-burbuja.jpg",
    "macetero-ceramica-gris.jpg",
    "canasto-ovalado.jpg"
)

for ($i = 0; $i -lt $fotos.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $fotos[$i]
}

# Column widths: new Foto column is wider; restore the widths of the
# shifted columns to their original (pre-shift) values. (The engine's
# ColumnWidth setter stores width + 5/6 internally, so subtract that
# offset to land on the clean target widths.)
$widthOffset = 0.8333333333333333
$ws.Columns.Item(7).ColumnWidth = 25 - $widthOffset
$ws.Columns.Item(8).ColumnWidth = 10 - $widthOffset
$ws.Columns.Item(9).ColumnWidth = 8 - $widthOffset
$ws.Columns.Item(10).ColumnWidth = 12 - $widthOffset
$ws.Columns.Item(11).ColumnWidth = 18 - $widthOffset
